$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = 'Datos actualizados a 31 de Mayo de 2020 a las 18:05'
$ws.Cells.Item(4, 2).Value = 1821162
$ws.Cells.Item(4, 3).Value = 4342
$ws.Cells.Item(4, 4).Value = 535387
$ws.Cells.Item(4, 5).Value = 1180121
$ws.Cells.Item(4, 7).Value = 97
$ws.Cells.Item(4, 8).Value = 105654
$ws.Cells.Item(10, 1).Value = 'India'
$ws.Cells.Item(10, 2).Value = 189717
$ws.Cells.Item(10, 3).Value = 7890
$ws.Cells.Item(10, 4).Value = 91016
$ws.Cells.Item(10, 5).Value = 93311
$ws.Cells.Item(10, 7).Value = 205
$ws.Cells.Item(10, 8).Value = 5390
$ws.Cells.Item(11, 1).Value = 'Francia'
$ws.Cells.Item(11, 2).Value = 188625
$ws.Cells.Item(11, 3).Value = 0
$ws.Cells.Item(11, 4).Value = 68268
$ws.Cells.Item(11, 5).Value = 91586
$ws.Cells.Item(11, 7).Value = 0
$ws.Cells.Item(11, 8).Value = 28771
$ws.Cells.Item(12, 2).Value = 183411
$ws.Cells.Item(12, 3).Value = 117
$ws.Cells.Item(12, 5).Value = 9609
$ws.Cells.Item(38, 2).Value = 23786
$ws.Cells.Item(38, 3).Value = 215
$ws.Cells.Item(38, 5).Value = 11451
$ws.Cells.Item(43, 4).Value = 10559
$ws.Cells.Item(43, 5).Value = 6224
$ws.Cells.Item(58, 1).Value = 'Argelia'
$ws.Cells.Item(58, 2).Value = 9394
$ws.Cells.Item(58, 3).Value = 127
$ws.Cells.Item(58, 4).Value = 5748
$ws.Cells.Item(58, 5).Value = 2993
$ws.Cells.Item(58, 7).Value = 7
$ws.Cells.Item(58, 8).Value = 653
$ws.Cells.Item(59, 1).Value = 'Armenia'
$ws.Cells.Item(59, 2).Value = 9282
$ws.Cells.Item(59, 3).Value = 355
$ws.Cells.Item(59, 4).Value = 3386
$ws.Cells.Item(59, 5).Value = 5765
$ws.Cells.Item(59, 7).Value = 4
$ws.Cells.Item(59, 8).Value = 131
$ws.Cells.Item(60, 2).Value = 9257
$ws.Cells.Item(60, 3).Value = 27
$ws.Cells.Item(60, 4).Value = 6551
$ws.Cells.Item(60, 5).Value = 2386
$ws.Cells.Item(60, 7).Value = 1
$ws.Cells.Item(60, 8).Value = 320
$ws.Cells.Item(68, 2).Value = 6439
$ws.Cells.Item(68, 3).Value = 260
$ws.Cells.Item(68, 4).Value = 3156
$ws.Cells.Item(68, 5).Value = 3078
$ws.Cells.Item(68, 7).Value = 10
$ws.Cells.Item(68, 8).Value = 205
$ws.Cells.Item(74, 2).Value = 4018
$ws.Cells.Item(74, 3).Value = 2
$ws.Cells.Item(74, 4).Value = 3833
$ws.Cells.Item(74, 5).Value = 75
$ws.Cells.Item(80, 2).Value = 3354
$ws.Cells.Item(80, 3).Value = 160
$ws.Cells.Item(80, 4).Value = 1504
$ws.Cells.Item(80, 5).Value = 1826
$ws.Cells.Item(80, 7).Value = 2
$ws.Cells.Item(80, 8).Value = 24
$ws.Cells.Item(91, 2).Value = 2045
$ws.Cells.Item(91, 3).Value = 20
$ws.Cells.Item(91, 4).Value = 1809
$ws.Cells.Item(91, 5).Value = 153
$ws.Cells.Item(92, 1).Value = 'Somalia'
$ws.Cells.Item(92, 2).Value = 1976
$ws.Cells.Item(92, 3).Value = 60
$ws.Cells.Item(92, 4).Value = 348
$ws.Cells.Item(92, 5).Value = 1550
$ws.Cells.Item(92, 7).Value = 5
$ws.Cells.Item(92, 8).Value = 78
$ws.Cells.Item(93, 1).Value = 'Kenia'
$ws.Cells.Item(93, 2).Value = 1962
$ws.Cells.Item(93, 3).Value = 74
$ws.Cells.Item(93, 4).Value = 478
$ws.Cells.Item(93, 5).Value = 1420
$ws.Cells.Item(93, 7).Value = 1
$ws.Cells.Item(93, 8).Value = 64
$ws.Cells.Item(101, 2).Value = 1631
$ws.Cells.Item(101, 3).Value = 18
$ws.Cells.Item(101, 5).Value = 820
$ws.Cells.Item(108, 1).Value = 'Mali'
$ws.Cells.Item(108, 2).Value = 1265
$ws.Cells.Item(108, 3).Value = 15
$ws.Cells.Item(108, 4).Value = 716
$ws.Cells.Item(108, 5).Value = 472
$ws.Cells.Item(108, 7).Value = 1
$ws.Cells.Item(108, 8).Value = 77
$ws.Cells.Item(109, 1).Value = 'Guinea-Bisau'
$ws.Cells.Item(109, 2).Value = 1256
$ws.Cells.Item(109, 4).Value = 42
$ws.Cells.Item(109, 5).Value = 1206
$ws.Cells.Item(109, 8).Value = 8
$ws.Cells.Item(118, 1).Value = 'Republica de Africa Central'
$ws.Cells.Item(118, 2).Value = 1011
$ws.Cells.Item(118, 3).Value = 49
$ws.Cells.Item(118, 4).Value = 23
$ws.Cells.Item(118, 5).Value = 986
$ws.Cells.Item(118, 7).Value = 1
$ws.Cells.Item(118, 8).Value = 2
$ws.Cells.Item(119, 1).Value = 'Sudan del Sur'
$ws.Cells.Item(119, 2).Value = 994
$ws.Cells.Item(119, 4).Value = 6
$ws.Cells.Item(119, 5).Value = 978
$ws.Cells.Item(119, 8).Value = 10
$ws.Cells.Item(120, 1).Value = 'Paraguay'
$ws.Cells.Item(120, 2).Value = 964
$ws.Cells.Item(120, 4).Value = 466
$ws.Cells.Item(120, 5).Value = 487
$ws.Cells.Item(120, 8).Value = 11
$ws.Cells.Item(122, 2).Value = 944
$ws.Cells.Item(122, 3).Value = 1
$ws.Cells.Item(122, 5).Value = 137
$ws.Cells.Item(123, 2).Value = 861
$ws.Cells.Item(123, 3).Value = 9
$ws.Cells.Item(123, 4).Value = 454
$ws.Cells.Item(123, 5).Value = 361
$ws.Cells.Item(131, 2).Value = 739
$ws.Cells.Item(131, 3).Value = 5
$ws.Cells.Item(131, 4).Value = 522
$ws.Cells.Item(131, 5).Value = 208
$ws.Cells.Item(200, 1).Value = 'Belice'
$ws.Cells.Item(200, 4).Value = 16
$ws.Cells.Item(200, 8).Value = 2
$ws.Cells.Item(201, 1).Value = 'Santa Lucia'
$ws.Cells.Item(201, 4).Value = 18
$ws.Cells.Item(201, 8).Value = 0
$ws.Cells.Item(213, 1).Value = 'Islas Virgenes Britanicas'
$ws.Cells.Item(213, 4).Value = 7
$ws.Cells.Item(213, 8).Value = 1
$ws.Cells.Item(214, 1).Value = 'Papua Nueva Guinea'
$ws.Cells.Item(214, 4).Value = 8
$ws.Cells.Item(214, 8).Value = 0
